$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 11 (network/patricia) label: runme_small.sh -> runme_large.sh
$ws.Range("A11").Value = "network/patricia/runme_large.sh"

# Add new row 12 for office/rsynth
$ws.Range("A12").Value = "office/rsynth/runme_large.sh"
$ws.Range("B12").Value = 0.13
$ws.Range("C12").Value = 0.14
$ws.Range("D12").Value = 0

# Update active selection to A21
$ws.Range("A21").Select()
